$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new value, derived from the update to this week's
# FlashScore odds/correct-score data (rows 2-22).
$updates = @{
    "AB2" = 12;
    "AC2" = 9.5;
    "AD2" = 19;
    "AE2" = 15;
    "AF2" = 21;
    "AG2" = 17;
    "AH2" = 8.5;
    "AL2" = 17;
    "AM2" = 23;
    "AP2" = 29;
    "G2" = 1.95;
    "I2" = 3.6;
    "J2" = 2.5;
    "L2" = 4;
    "O2" = 1.17;
    "P2" = 5;
    "Q2" = 1.57;
    "R2" = 2.38;
    "S2" = 1.95;
    "T2" = 1.95;
    "U2" = 2.38;
    "V2" = 1.57;
    "W2" = 1.29;
    "X2" = 3.5;
    "AA3" = 10;
    "AB3" = 13;
    "AD3" = 23;
    "AG3" = 12;
    "AH3" = 7.5;
    "AI3" = 15;
    "AK3" = 201;
    "AL3" = 12;
    "AO3" = 41;
    "AP3" = 26;
    "AQ3" = 34;
    "G3" = 2.2;
    "H3" = 3.5;
    "I3" = 3.1;
    "K3" = 2.2;
    "L3" = 3.6;
    "O3" = 1.25;
    "P3" = 4;
    "U3" = 3;
    "V3" = 1.4;
    "W3" = 1.36;
    "X3" = 3;
    "Y3" = 1.67;
    "Z3" = 2.1;
    "AA4" = 8.5;
    "AB4" = 13;
    "AD4" = 26;
    "AE4" = 21;
    "AG4" = 9.5;
    "AH4" = 7;
    "AL4" = 10;
    "AM4" = 17;
    "AN4" = 12;
    "AO4" = 34;
    "AP4" = 29;
    "AQ4" = 41;
    "G4" = 2.35;
    "I4" = 2.9;
    "J4" = 3.1;
    "L4" = 3.6;
    "O4" = 1.36;
    "P4" = 3.2;
    "U4" = 4;
    "V4" = 1.25;
    "AB5" = 9;
    "AC5" = 11;
    "AG5" = 26;
    "AH5" = 17;
    "AJ5" = 51;
    "AM5" = 67;
    "AN5" = 34;
    "AO5" = 151;
    "AQ5" = 51;
    "G5" = 1.18;
    "H5" = 7.5;
    "I5" = 13;
    "J5" = 1.53;
    "K5" = 3.2;
    "L5" = 9;
    "Q5" = 1.3;
    "R5" = 3.5;
    "U5" = 1.8;
    "V5" = 2;
    "W5" = 1.18;
    "X5" = 4.5;
    "AA7" = 5;
    "AC7" = 9;
    "AG7" = 7;
    "AI7" = 23;
    "AN7" = 21;
    "AP7" = 51;
    "AR7" = 1.78;
    "AS7" = 2.1;
    "H7" = 3.6;
    "J7" = 2.2;
    "K7" = 2.05;
    "L7" = 7;
    "Q7" = 2.35;
    "R7" = 1.57;
    "S7" = 3.65;
    "T7" = 1.27;
    "U7" = 4.33;
    "V7" = 1.2;
    "W7" = 1.5;
    "X7" = 2.5;
    "Y7" = 2.38;
    "Z7" = 1.53;
    "AC8" = 10;
    "AD8" = 26;
    "AE8" = 21;
    "AG8" = 9;
    "AM8" = 13;
    "AN8" = 11;
    "AO8" = 29;
    "AP8" = 23;
    "G8" = 2.63;
    "I8" = 2.75;
    "J8" = 3.25;
    "K8" = 2.1;
    "L8" = 3.4;
    "AC9" = 12;
    "AE9" = 26;
    "AG9" = 11;
    "AH9" = 6.5;
    "AK9" = 201;
    "AL9" = 8.5;
    "AQ9" = 26;
    "G9" = 3.2;
    "I9" = 2.15;
    "J9" = 3.75;
    "O9" = 1.29;
    "P9" = 3.75;
    "Q9" = 1.88;
    "R9" = 1.98;
    "U9" = 3.2;
    "V9" = 1.36;
    "Y9" = 1.67;
    "Z9" = 2.1;
    "M10" = 1.06;
    "N10" = 10;
    "Q10" = 2.1;
    "R10" = 1.7;
    "U10" = 3.75;
    "V10" = 1.25;
    "AL11" = 8;
    "Q11" = 1.93;
    "R11" = 1.93;
    "U11" = 3.25;
    "V11" = 1.33;
    "AA12" = 7.5;
    "AC12" = 11;
    "AD12" = 29;
    "AE12" = 26;
    "AG12" = 7.5;
    "AI12" = 17;
    "AK12" = 451;
    "AL12" = 7;
    "AM12" = 11;
    "AO12" = 26;
    "AR12" = 1.8;
    "AS12" = 2.05;
    "G12" = 2.8;
    "H12" = 3.2;
    "I12" = 2.5;
    "J12" = 3.75;
    "K12" = 1.95;
    "L12" = 3.4;
    "M12" = 1.08;
    "N12" = 8;
    "O12" = 1.44;
    "P12" = 2.63;
    "Q12" = 2.35;
    "R12" = 1.57;
    "U12" = 4.5;
    "V12" = 1.18;
    "W12" = 1.53;
    "X12" = 2.38;
    "Y12" = 2;
    "Z12" = 1.75;
    "AE14" = 12;
    "AL14" = 19;
    "G14" = 1.44;
    "J14" = 1.95;
    "K14" = 2.4;
    "M14" = 1.03;
    "N14" = 15;
    "Q14" = 1.67;
    "R14" = 2.15;
    "W14" = 1.33;
    "X14" = 3.25;
    "M18" = 1.04;
    "N18" = 9;
    "AD19" = 7;
    "AE19" = 11;
    "AF19" = 29;
    "AH19" = 19;
    "AI19" = 34;
    "AJ19" = 81;
    "AK19" = 401;
    "AL19" = 34;
    "G19" = 1.13;
    "I19" = 15;
    "L19" = 12;
    "N19" = 26;
    "U19" = 1.8;
    "V19" = 2;
    "Y19" = 2.05;
    "Z19" = 1.7;
    "AA22" = 6.5;
    "AB22" = 6.5;
    "AC22" = 9.5;
    "AD22" = 9.5;
    "AE22" = 13;
    "AF22" = 29;
    "AG22" = 11;
    "AH22" = 9;
    "AI22" = 21;
    "AJ22" = 67;
    "AK22" = 1000;
    "AL22" = 15;
    "AM22" = 34;
    "AN22" = 19;
    "AO22" = 67;
    "G22" = 1.45;
    "H22" = 4.5;
    "J22" = 1.95;
    "K22" = 2.25;
    "M22" = 1.04;
    "N22" = 9;
    "O22" = 1.25;
    "P22" = 3.75;
    "Q22" = 1.85;
    "R22" = 2;
    "U22" = 3;
    "V22" = 1.36;
    "W22" = 1.36;
    "X22" = 3;
    "Y22" = 2.1;
    "Z22" = 1.67;
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
